# Apply updated DUTResistance / Resistance / SingleDeviceResistance values
# (columns J, AB, AC) for rows 2-28, per the "speedup with lookup tables pt2"
# change. SingleDeviceResistance (AC) is now set equal to Resistance (AB/J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 27.99187006961076
    3  = 27.60218805472812
    4  = 27.54465938100303
    5  = 25.81421241681556
    6  = 25.38710696946947
    7  = 25.30950471977179
    8  = 24.82388561306556
    9  = 24.37519151106132
    10 = 24.28592585215918
    11 = 19.74970716676716
    12 = 19.31177628505161
    13 = 19.22833286994427
    14 = 17.32040946817778
    15 = 16.81784124483682
    16 = 16.69922881674465
    17 = 16.18113487099227
    18 = 15.63892530632404
    19 = 15.49851661728377
    20 = 17.64592457541496
    21 = 17.19137496712583
    22 = 17.09893612708201
    23 = 15.12814465154436
    24 = 14.59819771446141
    25 = 14.4645490643557
    26 = 13.93447080663309
    27 = 13.35737332032622
    28 = 13.19765080018907
}

foreach ($row in $newValues.Keys) {
    $value = $newValues[$row]
    $ws.Range("J$row").Value = $value
    $ws.Range("AB$row").Value = $value
    $ws.Range("AC$row").Value = $value
}
